# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
# Applies the league-data refresh to the "Australia ALeague" sheet:
#  - rows 73/74, 104/105, 124/125 swap their B..AC content (id in col A stays put)
#  - rows 154-157 get new match data (B, E-V columns)
#  - two brand-new rows (158, 159) are appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $r1, $r2, $firstCol, $lastCol) {
    # Read everything first so the swap doesn't clobber itself.
    $vals1 = @{}
    $vals2 = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals1[$c] = $ws.Cells.Item($r1, $c).Value2()
        $vals2[$c] = $ws.Cells.Item($r2, $c).Value2()
    }
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r1, $c).Value = $vals2[$c]
        $ws.Cells.Item($r2, $c).Value = $vals1[$c]
    }
}

# Columns B..AC are 2..29
Swap-RowData $ws 73 74 2 29
Swap-RowData $ws 104 105 2 29
Swap-RowData $ws 124 125 2 29

function Set-MatchRow {
    param(
        $ws,
        [int]$row,
        [long]$id,
        [double]$date,
        [string]$home,
        [string]$away,
        [double]$K, [double]$L, [double]$M, [double]$N, [double]$O, [double]$P,
        [double]$Q, [double]$R, [double]$S, [double]$T, [double]$U, [double]$V
    )
    $ws.Cells.Item($row, 2).Value = $id          # B id
    $ws.Cells.Item($row, 5).Value = $date        # E Date
    $ws.Cells.Item($row, 6).Value = $home        # F HomeTeam
    $ws.Cells.Item($row, 7).Value = $away        # G AwayTeam
    $ws.Cells.Item($row, 11).Value = $K          # K oddH_op
    $ws.Cells.Item($row, 12).Value = $L          # L oddD_op
    $ws.Cells.Item($row, 13).Value = $M          # M oddA_op
    $ws.Cells.Item($row, 14).Value = $N          # N oddH
    $ws.Cells.Item($row, 15).Value = $O          # O oddD
    $ws.Cells.Item($row, 16).Value = $P          # P oddA
    $ws.Cells.Item($row, 17).Value = $Q          # Q Ah
    $ws.Cells.Item($row, 18).Value = $R          # R oddAHH
    $ws.Cells.Item($row, 19).Value = $S          # S oddAHA
    $ws.Cells.Item($row, 20).Value = $T          # T AhOU
    $ws.Cells.Item($row, 21).Value = $U          # U oddAHOver
    $ws.Cells.Item($row, 22).Value = $V          # V oddAHUnder
}

# --- rows 154-157: existing rows get new fixtures/odds (B, E..V). ---
# H/I/J (result cols) and W..AA (PL cols) are untouched/absent for these, same as before.
Set-MatchRow $ws 154 7127416 45408.28125             "Brisbane Roar"       "Adelaide United"          2.3   4     2.5   2.1   4     2.75  -0.25 1.97  1.93  3.75  2.025 1.825
Set-MatchRow $ws 155 7127418 45409.17708333334        "Newcastle Jets"      "Central Coast Mariners"   3.6   3.25  2     3.8   3.3   1.909 0.5   1.93  1.97  3     1.925 1.925
Set-MatchRow $ws 156 7127419 45409.17708333334        "Wellington Phoenix"  "Macarthur FC"             1.85  3.5   3.9   1.75  3.6   4.2   -0.75 2.03  1.87  3     1.875 1.975
Set-MatchRow $ws 157 7127417 45409.28125             "Melbourne Victory"   "Western Sydney Wanderers" 2.05  3.3   3.4   1.833 3.5   4     -0.5  1.86  2.04  3     1.925 1.925

# --- rows 158-159: brand-new fixtures appended at the bottom of the table. ---
# Clone formatting (bold/bordered id cell, date-formatted date cell) from the row above,
# then overwrite with the real values.
$ws.Range("A157").Copy($ws.Range("A158"))
$ws.Range("E157").Copy($ws.Range("E158"))
$ws.Range("A157").Copy($ws.Range("A159"))
$ws.Range("E157").Copy($ws.Range("E159"))

$ws.Cells.Item(158, 1).Value = 156
$ws.Cells.Item(158, 3).Value = "Australia ALeague"
$ws.Cells.Item(158, 4).Value = "Australia ALeague"
Set-MatchRow $ws 158 8109525 45410.08333333334 "Sydney FC" "Perth Glory" 1.5 3.6 7 1.444 4 7.5 -1.25 1.9 2 3.5 1.875 1.975
$ws.Cells.Item(158, 23).Value = 0   # W
$ws.Cells.Item(158, 24).Value = 0   # X
$ws.Cells.Item(158, 25).Value = 0   # Y
$ws.Cells.Item(158, 26).Value = 0   # Z
$ws.Cells.Item(158, 27).Value = 0   # AA

$ws.Cells.Item(159, 1).Value = 157
$ws.Cells.Item(159, 3).Value = "Australia ALeague"
$ws.Cells.Item(159, 4).Value = "Australia ALeague"
Set-MatchRow $ws 159 7127421 45410.16666666666 "Melbourne City" "Western United FC" 1.65 4 4.333 1.4 4.2 7 -1.25 1.89 2.01 3.5 2.025 1.825
$ws.Cells.Item(159, 23).Value = 0   # W
$ws.Cells.Item(159, 24).Value = 0   # X
$ws.Cells.Item(159, 25).Value = 0   # Y
$ws.Cells.Item(159, 26).Value = 0   # Z
$ws.Cells.Item(159, 27).Value = 0   # AA

Write-Output "done"
